# Add the new "strategy_id-6004" worksheet as the last sheet in the workbook,
# populate it with the same layout/data as the other strategy sheets, and
# reuse the existing bold/bordered header style (style index 1) by copying
# formats from an existing sheet's header row instead of building new styles.

$wb = $excel.ActiveWorkbook

$templateSheet = $wb.Worksheets.Item(2)
$firstSheet    = $wb.Worksheets.Item(1)
$lastSheet     = $wb.Worksheets.Item($wb.Worksheets.Count)

$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "strategy_id-6004"

# ---- Header row (row 1) ----
$headers = @(
    "subsector",
    "variable",
    "variable_trajectory_group",
    "normalize_group",
    "trajgroup_no_vary_q",
    "uniform_scaling_q",
    "variable_trajectory_group_trajectory_type",
    "max_55",
    "min_55"
)
for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}
for ($n = 0; $n -le 55; $n++) {
    $ws.Cells.Item(1, 10 + $n).Value = $n
}

# ---- Data row (row 2) ----
$ws.Cells.Item(2, 1).Value = "General"
$ws.Cells.Item(2, 2).Value = "frac_gnrl_eating_red_meat"
$ws.Cells.Item(2, 3).Value = 13
$ws.Cells.Item(2, 4).Value = ""
$ws.Cells.Item(2, 5).Value = ""
$ws.Cells.Item(2, 6).Value = ""
$ws.Cells.Item(2, 7).Value = ""
$ws.Cells.Item(2, 8).Value = 1
$ws.Cells.Item(2, 9).Value = 1

$values55 = @(
    0.31,0.31,0.31,0.31,0.31,0.31,0.31,0.31,0.31,0.31,
    0.31,0.31,0.31,0.31,0.31,0.31,0.31,0.31,0.31,0.31,
    0.31,0.31,0.3100000000000001,0.31,0.31,0.3100000000000001,0.3099999999999999,0.31,0.31,0.3099999999999999,
    0.31,0.31,0.31,0.31,0.31,0.31,0.31,0.31,0.31,0.31,
    0.31,0.31,0.31,0.31,0.31,0.31,0.31,0.31,0.31,0.31,
    0.31,0.31,0.31,0.31,0.31
)
for ($i = 0; $i -lt $values55.Count; $i++) {
    $ws.Cells.Item(2, 10 + $i).Value = $values55[$i]
}

# ---- Formatting: reuse the bold/centered/bordered header style already
# used by the other strategy sheets (style index 1), by copying it from the
# template sheet's header row rather than re-declaring new style entries. ----
$templateSheet.Range("A1:BM1").Copy() | Out-Null
$ws.Range("A1:BM1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Restore the originally active sheet/selection so the workbook view state
# (active tab) is left exactly as before the edit.
$firstSheet.Activate()
$firstSheet.Range("A1").Select() | Out-Null
